# Apply stock-report corrections to CryCompanywiseStockReport.
# Columns: A=Sl.No, B=ItemCode, C=ItemName, D=Rate, E=MRP/Rate2, F=Qty, G=Value
# Row types: item rows, "Sub Total:" rows (value in column B), "Grand Total:" row (column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 3M INDIA LTD ---
# Row 5: 3M-SB FOOTLOCK MOP
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 4521.72
# Row 10: Sub Total
$ws.Range("B10").Value = 25049.08

# --- (IFB section) ---
# Row 18: IFB-Top Load Fully Automatic 7 Kg 5 Star TL - 701 A Series
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 29432.66
# Row 19: Sub Total
$ws.Range("B19").Value = 450860.53

# --- BUT-Kromo Deluxe KCp2 section ---
# Row 106
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 4800.05
# Row 115: Sub Total
$ws.Range("B115").Value = 79848.45

# --- CAV-MEERA HERBAL POWDER 150GM CONT section ---
# Row 119
$ws.Range("F119").Value = 226
$ws.Range("G119").Value = 11765.56
# Row 125: Sub Total
$ws.Range("B125").Value = 55484.48

# --- DAB-Real Activ Coconut Water Tetra 1000ml (rows 167/168 swap item codes) ---
# Row 167 (Sl.No 129)
$ws.Range("B167").Value = 64329
$ws.Range("E167").Value = 128.32
$ws.Range("F167").Value = 1
$ws.Range("G167").Value = 120.69
# Row 168 (Sl.No 130)
$ws.Range("B168").Value = 57552
$ws.Range("E168").Value = 136.86
$ws.Range("F168").Value = -5
$ws.Range("G168").Value = -603.45

# --- EVE-Eveready batteries ---
# Row 184: EVE-Eveready AA Battery 1015
$ws.Range("F184").Value = 850
$ws.Range("G184").Value = 7182.5
# Row 185: EVE-Eveready AAA Battery 1012
$ws.Range("F185").Value = 334
$ws.Range("G185").Value = 2668.66
# Row 187: Sub Total
$ws.Range("B187").Value = 10777.25

# --- HAM-Storex Container 1000 Ml ---
# Row 225
$ws.Range("F225").Value = 0
$ws.Range("G225").Value = 0
# Row 231: Sub Total
$ws.Range("B231").Value = 25144.39

# --- HIM-COCOA BUTTER INTEN.BODY LOTION 200ML ---
# Row 246
$ws.Range("F246").Value = 0
$ws.Range("G246").Value = 0

# --- HIM-LIP BALM (12S BLISTER PACK ) 10G ---
# Row 259
$ws.Range("F259").Value = 22
$ws.Range("G259").Value = 443.52
# Row 272: Sub Total
$ws.Range("B272").Value = 28073.18

# --- HUL-Kissan nango jam 490g (rows 283/284 swap item codes) ---
# Row 283 (Sl.No 225)
$ws.Range("B283").Value = 63520
$ws.Range("E283").Value = 153.4
$ws.Range("F283").Value = 63
$ws.Range("G283").Value = 9089.639999999999
# Row 284 (Sl.No 226)
$ws.Range("B284").Value = 55373
$ws.Range("E284").Value = 163.62
$ws.Range("F284").Value = -94
$ws.Range("G284").Value = -13562.32

# --- HUL-Knorr Schezwan Sauce 200Gm Pauch ---
# Row 290
$ws.Range("F290").Value = 29
$ws.Range("G290").Value = 1488.28
# Row 316: Sub Total
$ws.Range("B316").Value = 152882.66

# --- MAYA-AIR Gold 3B GT Gas Stove ---
# Row 430
$ws.Range("F430").Value = 2
$ws.Range("G430").Value = 10927.1

# --- MAYA-Table Top Grinder Jewel Stone ---
# Row 438
$ws.Range("F438").Value = 4
$ws.Range("G438").Value = 21237.4
# Row 445: Sub Total
$ws.Range("B445").Value = 247754.94

# --- CRE-Butter cookies 100gm (rows 448/449 swap item codes) ---
# Row 448 (Sl.No 364)
$ws.Range("B448").Value = 65068
$ws.Range("E448").Value = 13.97
$ws.Range("F448").Value = 63
$ws.Range("G448").Value = 828.45
# Row 449 (Sl.No 365)
$ws.Range("B449").Value = 53602
$ws.Range("E449").Value = 15.69
$ws.Range("F449").Value = -231
$ws.Range("G449").Value = -3037.65

# --- CRE-Cremica Pista Almond Cookies (75 +25Gm) (rows 466/467 swap item codes) ---
# Row 466 (Sl.No 382)
$ws.Range("B466").Value = 64919
$ws.Range("E466").Value = 27.97
$ws.Range("F466").Value = 61
$ws.Range("G466").Value = 1604.3
# Row 467 (Sl.No 383)
$ws.Range("B467").Value = 45702
$ws.Range("E467").Value = 31.43
$ws.Range("F467").Value = -215
$ws.Range("G467").Value = -5654.5

# --- CRE-Kaju khz cookies 100 gm (rows 469/470 swap item codes) ---
# Row 469 (Sl.No 385)
$ws.Range("B469").Value = 53595
$ws.Range("E469").Value = 17.61
$ws.Range("F469").Value = -335
$ws.Range("G469").Value = -4934.55
# Row 470 (Sl.No 386)
$ws.Range("B470").Value = 65067
$ws.Range("E470").Value = 15.65
$ws.Range("F470").Value = 126
$ws.Range("G470").Value = 1855.98

# --- OCT-Octavius Instant Coffee Gold 100 gm.(Granulated) ---
# Row 506
$ws.Range("F506").Value = 6
$ws.Range("G506").Value = 1296.66
# Row 508: Sub Total
$ws.Range("B508").Value = 1729.33

# --- Rasna Nagpur Orange (32 Glass) (rows 573/574 swap item codes) ---
# Row 573 (Sl.No 465)
$ws.Range("B573").Value = 60022
$ws.Range("E573").Value = 37.22
$ws.Range("F573").Value = -113
$ws.Range("G573").Value = -3709.79
# Row 574 (Sl.No 466)
$ws.Range("B574").Value = 64830
$ws.Range("E574").Value = 34.9
$ws.Range("F574").Value = 101
$ws.Range("G574").Value = 3315.83

# --- RKT- Moov 80 gm Spray ---
# Row 585
$ws.Range("F585").Value = 4
$ws.Range("G585").Value = 701.88
# Row 588: Sub Total
$ws.Range("B588").Value = 27695.94

# --- SOU-Black Fard Dates 400g ---
# Row 626
$ws.Range("F626").Value = 540
$ws.Range("G626").Value = 65799

# --- SOU-Kodo Millets 500 g ---
# Row 629
$ws.Range("F629").Value = 87
$ws.Range("G629").Value = 7443.72

# --- SOU-Little Millets 500g ---
# Row 630
$ws.Range("F630").Value = 79
$ws.Range("G630").Value = 6548.31
# Row 633: Sub Total
$ws.Range("B633").Value = 119240.63

# --- Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm (rows 721/722 swap item codes) ---
# Row 721 (Sl.No 579)
$ws.Range("B721").Value = 65362
$ws.Range("F721").Value = 18
$ws.Range("G721").Value = 735.66
# Row 722 (Sl.No 580)
$ws.Range("B722").Value = 65079
$ws.Range("F722").Value = 6
$ws.Range("G722").Value = 245.22

# --- VVD Priyam Cold Pressed Groundnut Oil Pouch 1 Ltr ---
# Row 751
$ws.Range("F751").Value = 3489
$ws.Range("G751").Value = 569090.79
# Row 758: Sub Total
$ws.Range("B758").Value = 695539.7

# --- WPL-230 IMPRO PRM 3S-z Single door Refriegrator ---
# Row 761
$ws.Range("F761").Value = 2
$ws.Range("G761").Value = 26733.76
# Row 764: Sub Total
$ws.Range("B764").Value = 185345.88

# --- Grand totals ---
# Row 777: Sub Total
$ws.Range("B777").Value = 5427370.3
# Row 778: Grand Total
$ws.Range("B778").Value = 5427370.3
